$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.987341780766807
$ws.Range("D16").Value = 1.070916661296198
$ws.Range("E16").Value = 0.9739138526456652
$ws.Range("F16").Value = 0.987341780766807
$ws.Range("G16").Value = 1.033966248719508
$ws.Range("H16").Value = 0.9478199508561423
$ws.Range("I16").Value = 0.9738484867596405
$ws.Range("J16").Value = 1.070916661296198
$ws.Range("K16").Value = 1.022415256970931
$ws.Range("L16").Value = 1.004878518868869
$ws.Range("M16").Value = 0.9979678301739935

# Copy the style/format from A15 (which carries the bold/border/centered style) onto A16
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
